$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWebsiteQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $newWebsiteQuery

$ws.Range("B2").Select() | Out-Null
